# Refresh market-price-derived columns (H:N) on several Leve-profit sheets.
# Values below mirror an updated Universalis price pull; CUL is untouched.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (48 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 999  # H32  was 2499.5
$ws.Cells.Item(32, 9).Value = 999  # I32  was 1000
$ws.Cells.Item(32, 10).Value = 0  # J32  was 3999
$ws.Cells.Item(32, 11).Value = 999  # K32  was 1000
$ws.Cells.Item(32, 12).Value = 0  # L32  was 3999
$ws.Cells.Item(32, 13).Value = -673  # M32  was -674
$ws.Cells.Item(32, 14).ClearContents()  # N32  was -4651
$ws.Cells.Item(38, 8).Value = 8474.174000000001  # H38  was 9337.15
$ws.Cells.Item(38, 9).Value = 10982.091  # I38  was 13205.556
$ws.Cells.Item(38, 10).Value = 6175.25  # J38  was 6172.091
$ws.Cells.Item(38, 11).Value = 32946.273  # K38  was 39616.66800000001
$ws.Cells.Item(38, 12).Value = 18525.75  # L38  was 18516.273
$ws.Cells.Item(38, 13).Value = -32574.273  # M38  was -39244.66800000001
$ws.Cells.Item(38, 14).Value = -19269.75  # N38  was -19260.273
$ws.Cells.Item(40, 8).Value = 4380.1763  # H40  was 4359.0557
$ws.Cells.Item(40, 10).Value = 4676.643  # J40  was 4631.533
$ws.Cells.Item(40, 12).Value = 4676.643  # L40  was 4631.533
$ws.Cells.Item(40, 14).Value = -5026.643  # N40  was -4981.533
$ws.Cells.Item(51, 8).Value = 3734.6155  # H51  was 3675
$ws.Cells.Item(51, 10).Value = 4660  # J51  was 4366.6665
$ws.Cells.Item(51, 12).Value = 4660  # L51  was 4366.6665
$ws.Cells.Item(51, 14).Value = -5628  # N51  was -5334.6665
$ws.Cells.Item(113, 8).Value = 4751.25  # H113  was 5333.3335
$ws.Cells.Item(113, 9).Value = 3005  # I113  was 0
$ws.Cells.Item(113, 11).Value = 3005  # K113  was 0
$ws.Cells.Item(113, 13).Value = 249  # M113  was (empty)
$ws.Cells.Item(116, 8).Value = 6047.25  # H116  was 5772.875
$ws.Cells.Item(116, 9).Value = 5396.3335  # I116  was 5531.5
$ws.Cells.Item(116, 10).Value = 8000  # J116  was 6497
$ws.Cells.Item(116, 11).Value = 5396.3335  # K116  was 5531.5
$ws.Cells.Item(116, 12).Value = 8000  # L116  was 6497
$ws.Cells.Item(116, 13).Value = -1954.3335  # M116  was -2089.5
$ws.Cells.Item(116, 14).Value = -14884  # N116  was -13381
$ws.Cells.Item(125, 8).Value = 103500  # H125  was 200000
$ws.Cells.Item(125, 10).Value = 103500  # J125  was 200000
$ws.Cells.Item(125, 12).Value = 931500  # L125  was 1800000
$ws.Cells.Item(125, 14).Value = -936420  # N125  was -1804920
$ws.Cells.Item(137, 8).Value = 1433.3793  # H137  was 1405.1333
$ws.Cells.Item(137, 9).Value = 1513.826  # I137  was 1475.125
$ws.Cells.Item(137, 10).Value = 1125  # J137  was 1125.1666
$ws.Cells.Item(137, 11).Value = 4541.478  # K137  was 4425.375
$ws.Cells.Item(137, 12).Value = 3375  # L137  was 3375.4998
$ws.Cells.Item(137, 13).Value = -1991.478  # M137  was -1875.375
$ws.Cells.Item(137, 14).Value = -8475  # N137  was -8475.4998
$ws.Cells.Item(141, 8).Value = 3366.1177  # H141  was 3267.7778
$ws.Cells.Item(141, 10).Value = 4716.3335  # J141  was 4270.5713
$ws.Cells.Item(141, 12).Value = 14149.0005  # L141  was 12811.7139
$ws.Cells.Item(141, 14).Value = -24509.0005  # N141  was -23171.7139

# --- Sheet: ARM (61 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2221.4849  # H32  was 2304.2344
$ws.Cells.Item(32, 9).Value = 2221.4849  # I32  was 2304.2344
$ws.Cells.Item(32, 11).Value = 2221.4849  # K32  was 2304.2344
$ws.Cells.Item(32, 13).Value = -1934.4849  # M32  was -2017.2344
$ws.Cells.Item(50, 8).Value = 5009.091  # H50  was 1922
$ws.Cells.Item(50, 9).Value = 809.4  # I50  was 694.5
$ws.Cells.Item(50, 10).Value = 8508.833000000001  # J50  was 3763.25
$ws.Cells.Item(50, 11).Value = 809.4  # K50  was 694.5
$ws.Cells.Item(50, 12).Value = 8508.833000000001  # L50  was 3763.25
$ws.Cells.Item(50, 13).Value = -95.39999999999998  # M50  was 19.5
$ws.Cells.Item(50, 14).Value = -9936.833000000001  # N50  was -5191.25
$ws.Cells.Item(61, 8).Value = 4277.579  # H61  was 4222.316
$ws.Cells.Item(61, 9).Value = 4414.1665  # I61  was 4222.316
$ws.Cells.Item(61, 10).Value = 1819  # J61  was 0
$ws.Cells.Item(61, 11).Value = 4414.1665  # K61  was 4222.316
$ws.Cells.Item(61, 12).Value = 1819  # L61  was 0
$ws.Cells.Item(61, 13).Value = -4202.1665  # M61  was -4010.316
$ws.Cells.Item(61, 14).Value = -2243  # N61  was (empty)
$ws.Cells.Item(74, 8).Value = 1990.25  # H74  was 2068.842
$ws.Cells.Item(74, 9).Value = 2123.875  # I74  was 2347.9285
$ws.Cells.Item(74, 10).Value = 1455.75  # J74  was 1287.4
$ws.Cells.Item(74, 11).Value = 2123.875  # K74  was 2347.9285
$ws.Cells.Item(74, 12).Value = 1455.75  # L74  was 1287.4
$ws.Cells.Item(74, 13).Value = -1249.875  # M74  was -1473.9285
$ws.Cells.Item(74, 14).Value = -3203.75  # N74  was -3035.4
$ws.Cells.Item(77, 8).Value = 1990.25  # H77  was 2068.842
$ws.Cells.Item(77, 9).Value = 2123.875  # I77  was 2347.9285
$ws.Cells.Item(77, 10).Value = 1455.75  # J77  was 1287.4
$ws.Cells.Item(77, 11).Value = 10619.375  # K77  was 11739.6425
$ws.Cells.Item(77, 12).Value = 7278.75  # L77  was 6437
$ws.Cells.Item(77, 13).Value = -6251.375  # M77  was -7371.6425
$ws.Cells.Item(77, 14).Value = -16014.75  # N77  was -15173
$ws.Cells.Item(102, 8).Value = 4059.6  # H102  was 4545.0454
$ws.Cells.Item(102, 9).Value = 2548.647  # I102  was 2987.7144
$ws.Cells.Item(102, 11).Value = 2548.647  # K102  was 2987.7144
$ws.Cells.Item(102, 13).Value = -926.6469999999999  # M102  was -1365.7144
$ws.Cells.Item(121, 8).Value = 54000  # H121  was 83418.336
$ws.Cells.Item(121, 9).Value = 54000  # I121  was 0
$ws.Cells.Item(121, 10).Value = 0  # J121  was 83418.336
$ws.Cells.Item(121, 11).Value = 54000  # K121  was 0
$ws.Cells.Item(121, 12).Value = 0  # L121  was 83418.336
$ws.Cells.Item(121, 13).Value = -52253  # M121  was (empty)
$ws.Cells.Item(121, 14).ClearContents()  # N121  was -86912.336
$ws.Cells.Item(122, 8).Value = 1469.037  # H122  was 1510.1538
$ws.Cells.Item(122, 9).Value = 1138.1364  # I122  was 1173.2858
$ws.Cells.Item(122, 11).Value = 3414.4092  # K122  was 3519.8574
$ws.Cells.Item(122, 13).Value = -964.4092000000001  # M122  was -1069.8574
$ws.Cells.Item(132, 8).Value = 1735.0646  # H132  was 1804.138
$ws.Cells.Item(132, 9).Value = 1725.3148  # I132  was 1781.1154
$ws.Cells.Item(132, 10).Value = 1800.875  # J132  was 2003.6666
$ws.Cells.Item(132, 11).Value = 5175.9444  # K132  was 5343.3462
$ws.Cells.Item(132, 12).Value = 5402.625  # L132  was 6010.9998
$ws.Cells.Item(132, 13).Value = -2645.9444  # M132  was -2813.3462
$ws.Cells.Item(132, 14).Value = -10462.625  # N132  was -11070.9998
$ws.Cells.Item(136, 8).Value = 4277.579  # H136  was 4222.316
$ws.Cells.Item(136, 9).Value = 4414.1665  # I136  was 4222.316
$ws.Cells.Item(136, 10).Value = 1819  # J136  was 0
$ws.Cells.Item(136, 11).Value = 13242.4995  # K136  was 12666.948
$ws.Cells.Item(136, 12).Value = 5457  # L136  was 0
$ws.Cells.Item(136, 13).Value = -10692.4995  # M136  was -10116.948
$ws.Cells.Item(136, 14).Value = -10557  # N136  was (empty)

# --- Sheet: BSM (15 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 29415394  # H86  was 33337180
$ws.Cells.Item(86, 9).Value = 3406.1428  # I86  was 3640.5833
$ws.Cells.Item(86, 11).Value = 3406.1428  # K86  was 3640.5833
$ws.Cells.Item(86, 13).Value = -2283.1428  # M86  was -2517.5833
$ws.Cells.Item(89, 8).Value = 29415394  # H89  was 33337180
$ws.Cells.Item(89, 9).Value = 3406.1428  # I89  was 3640.5833
$ws.Cells.Item(89, 11).Value = 17030.714  # K89  was 18202.9165
$ws.Cells.Item(89, 13).Value = -11414.714  # M89  was -12586.9165
$ws.Cells.Item(105, 8).Value = 4072.1667  # H105  was 3839.7144
$ws.Cells.Item(105, 9).Value = 3358.375  # I105  was 3306.4443
$ws.Cells.Item(105, 10).Value = 5499.75  # J105  was 4799.6
$ws.Cells.Item(105, 11).Value = 3358.375  # K105  was 3306.4443
$ws.Cells.Item(105, 12).Value = 5499.75  # L105  was 4799.6
$ws.Cells.Item(105, 13).Value = -1611.375  # M105  was -1559.4443
$ws.Cells.Item(105, 14).Value = -8993.75  # N105  was -8293.6

# --- Sheet: CRP (15 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(93, 8).Value = 37041.2  # H93  was 43699.75
$ws.Cells.Item(93, 9).Value = 33801.5  # I93  was 41599.668
$ws.Cells.Item(93, 11).Value = 33801.5  # K93  was 41599.668
$ws.Cells.Item(93, 13).Value = -31929.5  # M93  was -39727.668
$ws.Cells.Item(103, 8).Value = 12880.75  # H103  was 13666.333
$ws.Cells.Item(103, 9).Value = 12880.75  # I103  was 13666.333
$ws.Cells.Item(103, 11).Value = 12880.75  # K103  was 13666.333
$ws.Cells.Item(103, 13).Value = -11708.75  # M103  was -12494.333
$ws.Cells.Item(122, 8).Value = 951  # H122  was 1066.2727
$ws.Cells.Item(122, 9).Value = 818.1111  # I122  was 907.8
$ws.Cells.Item(122, 10).Value = 1100.5  # J122  was 1198.3334
$ws.Cells.Item(122, 11).Value = 2454.3333  # K122  was 2723.4
$ws.Cells.Item(122, 12).Value = 3301.5  # L122  was 3595.0002
$ws.Cells.Item(122, 13).Value = -4.333299999999781  # M122  was -273.3999999999996
$ws.Cells.Item(122, 14).Value = -8201.5  # N122  was -8495.0002

# --- Sheet: GSM (22 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 12466.167  # H102  was 13399.6
$ws.Cells.Item(102, 9).Value = 8699.5  # I102  was 8999.666999999999
$ws.Cells.Item(102, 11).Value = 8699.5  # K102  was 8999.666999999999
$ws.Cells.Item(102, 13).Value = -7077.5  # M102  was -7377.666999999999
$ws.Cells.Item(122, 8).Value = 24889.354  # H122  was 31192
$ws.Cells.Item(122, 9).Value = 29843.5  # I122  was 38444.223
$ws.Cells.Item(122, 10).Value = 12999.4  # J122  was 14874.5
$ws.Cells.Item(122, 11).Value = 89530.5  # K122  was 115332.669
$ws.Cells.Item(122, 12).Value = 38998.2  # L122  was 44623.5
$ws.Cells.Item(122, 13).Value = -87080.5  # M122  was -112882.669
$ws.Cells.Item(122, 14).Value = -43898.2  # N122  was -49523.5
$ws.Cells.Item(126, 8).Value = 3902.5  # H126  was 3100
$ws.Cells.Item(126, 9).Value = 3902.5  # I126  was 3100
$ws.Cells.Item(126, 11).Value = 11707.5  # K126  was 9300
$ws.Cells.Item(126, 13).Value = -9237.5  # M126  was -6830
$ws.Cells.Item(132, 8).Value = 2106.8206  # H132  was 2147.6943
$ws.Cells.Item(132, 9).Value = 2060.2222  # I132  was 2127
$ws.Cells.Item(132, 10).Value = 2666  # J132  was 2499.5
$ws.Cells.Item(132, 11).Value = 6180.6666  # K132  was 6381
$ws.Cells.Item(132, 12).Value = 7998  # L132  was 7498.5
$ws.Cells.Item(132, 13).Value = -3650.6666  # M132  was -3851
$ws.Cells.Item(132, 14).Value = -13058  # N132  was -12558.5

# --- Sheet: LTW (63 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 5780.4736  # H7  was 6342.6875
$ws.Cells.Item(7, 9).Value = 4169.4165  # I7  was 4631.8887
$ws.Cells.Item(7, 11).Value = 4169.4165  # K7  was 4631.8887
$ws.Cells.Item(7, 13).Value = -4057.4165  # M7  was -4519.8887
$ws.Cells.Item(10, 8).Value = 1503385.4  # H10  was 941044.25
$ws.Cells.Item(10, 9).Value = 3750075  # I10  was 1875526
$ws.Cells.Item(10, 10).Value = 5592.3335  # J10  was 6562.5
$ws.Cells.Item(10, 11).Value = 3750075  # K10  was 1875526
$ws.Cells.Item(10, 12).Value = 5592.3335  # L10  was 6562.5
$ws.Cells.Item(10, 13).Value = -3749935  # M10  was -1875386
$ws.Cells.Item(10, 14).Value = -5872.3335  # N10  was -6842.5
$ws.Cells.Item(20, 8).Value = 499.8  # H20  was 499.5
$ws.Cells.Item(20, 10).Value = 499.8  # J20  was 499.5
$ws.Cells.Item(20, 12).Value = 499.8  # L20  was 499.5
$ws.Cells.Item(20, 14).Value = -951.8  # N20  was -951.5
$ws.Cells.Item(22, 8).Value = 3343.9678  # H22  was 3569.074
$ws.Cells.Item(22, 9).Value = 3144.4614  # I22  was 3331.4167
$ws.Cells.Item(22, 10).Value = 3488.0557  # J22  was 3759.2
$ws.Cells.Item(22, 11).Value = 3144.4614  # K22  was 3331.4167
$ws.Cells.Item(22, 12).Value = 3488.0557  # L22  was 3759.2
$ws.Cells.Item(22, 13).Value = -2849.4614  # M22  was -3036.4167
$ws.Cells.Item(22, 14).Value = -4078.0557  # N22  was -4349.2
$ws.Cells.Item(27, 8).Value = 3343.9678  # H27  was 3569.074
$ws.Cells.Item(27, 9).Value = 3144.4614  # I27  was 3331.4167
$ws.Cells.Item(27, 10).Value = 3488.0557  # J27  was 3759.2
$ws.Cells.Item(27, 11).Value = 3144.4614  # K27  was 3331.4167
$ws.Cells.Item(27, 12).Value = 3488.0557  # L27  was 3759.2
$ws.Cells.Item(27, 13).Value = -3037.4614  # M27  was -3224.4167
$ws.Cells.Item(27, 14).Value = -3702.0557  # N27  was -3973.2
$ws.Cells.Item(55, 9).Value = 249  # I55  was 229.66667
$ws.Cells.Item(55, 10).Value = 3345  # J55  was 3999
$ws.Cells.Item(55, 11).Value = 249  # K55  was 229.66667
$ws.Cells.Item(55, 12).Value = 3345  # L55  was 3999
$ws.Cells.Item(55, 13).Value = -76  # M55  was -56.66667000000001
$ws.Cells.Item(55, 14).Value = -3691  # N55  was -4345
$ws.Cells.Item(61, 8).Value = 6171.5  # H61  was 6293.4
$ws.Cells.Item(82, 8).Value = 3393  # H82  was 3996.5833
$ws.Cells.Item(82, 9).Value = 2239.5  # I82  was 2870
$ws.Cells.Item(82, 10).Value = 4162  # J82  was 4559.875
$ws.Cells.Item(82, 11).Value = 2239.5  # K82  was 2870
$ws.Cells.Item(82, 12).Value = 4162  # L82  was 4559.875
$ws.Cells.Item(82, 13).Value = -1878.5  # M82  was -2509
$ws.Cells.Item(82, 14).Value = -4884  # N82  was -5281.875
$ws.Cells.Item(85, 8).Value = 3393  # H85  was 3996.5833
$ws.Cells.Item(85, 9).Value = 2239.5  # I85  was 2870
$ws.Cells.Item(85, 10).Value = 4162  # J85  was 4559.875
$ws.Cells.Item(85, 11).Value = 2239.5  # K85  was 2870
$ws.Cells.Item(85, 12).Value = 4162  # L85  was 4559.875
$ws.Cells.Item(85, 13).Value = -991.5  # M85  was -1622
$ws.Cells.Item(85, 14).Value = -6658  # N85  was -7055.875
$ws.Cells.Item(93, 8).Value = 3645.0833  # H93  was 3885.5908
$ws.Cells.Item(93, 9).Value = 1541.1538  # I93  was 1639.6364
$ws.Cells.Item(93, 11).Value = 1541.1538  # K93  was 1639.6364
$ws.Cells.Item(93, 13).Value = -293.1538  # M93  was -391.6364000000001
$ws.Cells.Item(113, 8).Value = 6171.5  # H113  was 6293.4
$ws.Cells.Item(126, 8).Value = 5780.4736  # H126  was 6342.6875
$ws.Cells.Item(126, 9).Value = 4169.4165  # I126  was 4631.8887
$ws.Cells.Item(126, 11).Value = 12508.2495  # K126  was 13895.6661
$ws.Cells.Item(126, 13).Value = -10038.2495  # M126  was -11425.6661
$ws.Cells.Item(141, 8).Value = 49999.5  # H141  was 0
$ws.Cells.Item(141, 10).Value = 49999.5  # J141  was 0
$ws.Cells.Item(141, 12).Value = 49999.5  # L141  was 0
$ws.Cells.Item(141, 14).Value = -60359.5  # N141  was (empty)

# --- Sheet: WVR (29 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 3446.9167  # H107  was 3279.4167
$ws.Cells.Item(107, 9).Value = 3285.5715  # I107  was 2998.625
$ws.Cells.Item(107, 10).Value = 3672.8  # J107  was 3841
$ws.Cells.Item(107, 11).Value = 9856.7145  # K107  was 8995.875
$ws.Cells.Item(107, 12).Value = 11018.4  # L107  was 11523
$ws.Cells.Item(107, 13).Value = -7936.7145  # M107  was -7075.875
$ws.Cells.Item(107, 14).Value = -14858.4  # N107  was -15363
$ws.Cells.Item(113, 8).Value = 1526.1111  # H113  was 1580.8823
$ws.Cells.Item(113, 10).Value = 1799.4445  # J113  was 1950
$ws.Cells.Item(113, 12).Value = 5398.333500000001  # L113  was 5850
$ws.Cells.Item(113, 14).Value = -9738.333500000001  # N113  was -10190
$ws.Cells.Item(122, 8).Value = 5922.5  # H122  was 6340
$ws.Cells.Item(122, 9).Value = 3719.375  # I122  was 3822.1428
$ws.Cells.Item(122, 10).Value = 8125.625  # J122  was 8857.857
$ws.Cells.Item(122, 11).Value = 11158.125  # K122  was 11466.4284
$ws.Cells.Item(122, 12).Value = 24376.875  # L122  was 26573.571
$ws.Cells.Item(122, 13).Value = -8708.125  # M122  was -9016.428400000001
$ws.Cells.Item(122, 14).Value = -29276.875  # N122  was -31473.571
$ws.Cells.Item(126, 8).Value = 1880.0714  # H126  was 2094.7273
$ws.Cells.Item(126, 9).Value = 1387.2  # I126  was 1513.2858
$ws.Cells.Item(126, 11).Value = 4161.6  # K126  was 4539.857400000001
$ws.Cells.Item(126, 13).Value = -1691.6  # M126  was -2069.857400000001
$ws.Cells.Item(136, 8).Value = 946.6923  # H136  was 985.26086
$ws.Cells.Item(136, 9).Value = 972.2273  # I136  was 1002.3
$ws.Cells.Item(136, 10).Value = 806.25  # J136  was 871.6667
$ws.Cells.Item(136, 11).Value = 2916.6819  # K136  was 3006.9
$ws.Cells.Item(136, 12).Value = 2418.75  # L136  was 2615.0001
$ws.Cells.Item(136, 13).Value = -366.6819  # M136  was -456.8999999999996
$ws.Cells.Item(136, 14).Value = -7518.75  # N136  was -7715.0001

# Total cells touched: 253
